# Applies commit: "Atualizacao de bases das ligas, do dia: 15-04-2024 as 22:35"
# - Rows 235-239 get reordered (re-sorted) while keeping the sequential id in column A.
# - Row 270 (previously a scheduled/un-played fixture) is updated with final match data.
# - Two brand-new fixture rows (271, 272) are appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 235-239: re-sort order (id in column A stays the same; B..AC columns permuted) ---
# Row 235
$ws.Range("B235").Value = 6852370
$ws.Range("C235").Value = "Romania Liga I"
$ws.Range("D235").Value = "Romania Liga I"
$ws.Range("E235").Value = 45359.625
$ws.Range("F235").Value = "Dinamo Bucharest"
$ws.Range("G235").Value = "ACS UTA Batrana Doamna"
$ws.Range("H235").Value = 1
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "H"
$ws.Range("K235").Value = 2.55
$ws.Range("L235").Value = 2.875
$ws.Range("M235").Value = 3
$ws.Range("N235").Value = 2.375
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = -0.25
$ws.Range("R235").Value = 2
$ws.Range("S235").Value = 1.85
$ws.Range("T235").Value = 2.25
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 1.875
$ws.Range("W235").Value = 1.375
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 1
$ws.Range("AA235").Value = -1
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 0.875

# Row 236
$ws.Range("B236").Value = 6836277
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = "Romania Liga I"
$ws.Range("E236").Value = 45359.625
$ws.Range("F236").Value = "CFR Cluj"
$ws.Range("G236").Value = "AFC Hermannstadt"
$ws.Range("H236").Value = 1
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 1.7
$ws.Range("L236").Value = 3.4
$ws.Range("M236").Value = 5
$ws.Range("N236").Value = 1.65
$ws.Range("O236").Value = 3.5
$ws.Range("P236").Value = 5.25
$ws.Range("Q236").Value = -0.75
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = 1.975
$ws.Range("W236").Value = 0.6499999999999999
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.425
$ws.Range("AA236").Value = -0.5
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.9750000000000001

# Row 237
$ws.Range("B237").Value = 6870268
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = "Romania Liga I"
$ws.Range("E237").Value = 45359.625
$ws.Range("F237").Value = "Petrolul Ploiesti"
$ws.Range("G237").Value = "ACS Sepsi"
$ws.Range("H237").Value = 1
$ws.Range("I237").Value = 2
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 2.8
$ws.Range("L237").Value = 3
$ws.Range("M237").Value = 2.55
$ws.Range("N237").Value = 3
$ws.Range("O237").Value = 3.2
$ws.Range("P237").Value = 2.3
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.85
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.875
$ws.Range("V237").Value = 1.975
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 1.3
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 1
$ws.Range("AB237").Value = 0.875
$ws.Range("AC237").Value = -1

# Row 238
$ws.Range("B238").Value = 6865915
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = "Romania Liga I"
$ws.Range("E238").Value = 45359.625
$ws.Range("F238").Value = "FC Voluntari"
$ws.Range("G238").Value = "Universitatea Cluj"
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 0
$ws.Range("J238").Value = "D"
$ws.Range("K238").Value = 3.5
$ws.Range("L238").Value = 3.25
$ws.Range("M238").Value = 2.05
$ws.Range("N238").Value = 3.4
$ws.Range("O238").Value = 3.1
$ws.Range("P238").Value = 2.15
$ws.Range("Q238").Value = 0.25
$ws.Range("R238").Value = 1.975
$ws.Range("S238").Value = 1.875
$ws.Range("T238").Value = 2.25
$ws.Range("U238").Value = 2.05
$ws.Range("V238").Value = 1.75
$ws.Range("W238").Value = -1
$ws.Range("X238").Value = 2.1
$ws.Range("Y238").Value = -1
$ws.Range("Z238").Value = 0.4875
$ws.Range("AA238").Value = -0.5
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 0.75

# Row 239
$ws.Range("B239").Value = 6861095
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = "Romania Liga I"
$ws.Range("E239").Value = 45359.625
$ws.Range("F239").Value = "FC Botosani"
$ws.Range("G239").Value = "Farul Constanta"
$ws.Range("H239").Value = 0
$ws.Range("I239").Value = 0
$ws.Range("J239").Value = "D"
$ws.Range("K239").Value = 3.75
$ws.Range("L239").Value = 3.4
$ws.Range("M239").Value = 1.909
$ws.Range("N239").Value = 3.1
$ws.Range("O239").Value = 3
$ws.Range("P239").Value = 2.375
$ws.Range("Q239").Value = 0.25
$ws.Range("R239").Value = 1.775
$ws.Range("S239").Value = 2.1
$ws.Range("T239").Value = 2
$ws.Range("U239").Value = 1.8
$ws.Range("V239").Value = 2.05
$ws.Range("W239").Value = -1
$ws.Range("X239").Value = 2
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = 0.3875
$ws.Range("AA239").Value = -0.5
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 1.05

# --- Row 270: update existing fixture now that the match has been played ---
# Row 270
$ws.Range("B270").Value = 7951791
$ws.Range("C270").Value = "Romania Liga I"
$ws.Range("D270").Value = "Romania Liga I"
$ws.Range("E270").Value = 45395.64583333334
$ws.Range("F270").Value = "Universitatea Cluj"
$ws.Range("G270").Value = "Petrolul Ploiesti"
$ws.Range("H270").Value = 1
$ws.Range("I270").Value = 2
$ws.Range("J270").Value = "A"
$ws.Range("K270").Value = 1.909
$ws.Range("L270").Value = 3.2
$ws.Range("M270").Value = 4
$ws.Range("N270").Value = 1.909
$ws.Range("O270").Value = 3.2
$ws.Range("P270").Value = 4.2
$ws.Range("Q270").Value = -0.5
$ws.Range("R270").Value = 1.925
$ws.Range("S270").Value = 1.925
$ws.Range("T270").Value = 2
$ws.Range("U270").Value = 1.95
$ws.Range("V270").Value = 1.9
$ws.Range("W270").Value = -1
$ws.Range("X270").Value = -1
$ws.Range("Y270").Value = 3.2
$ws.Range("Z270").Value = -1
$ws.Range("AA270").Value = 0.925
$ws.Range("AB270").Value = 0.95
$ws.Range("AC270").Value = -1

# --- Rows 271 and 272: brand-new match rows appended at the end ---
# Row 271
$ws.Range("A271").Value = 269
$ws.Range("B271").Value = 7951788
$ws.Range("C271").Value = "Romania Liga I"
$ws.Range("D271").Value = "Romania Liga I"
$ws.Range("E271").Value = 45396.3125
$ws.Range("F271").Value = "FC Botosani"
$ws.Range("G271").Value = "FC U Craiova 1948"
$ws.Range("H271").Value = 4
$ws.Range("I271").Value = 1
$ws.Range("J271").Value = "H"
$ws.Range("K271").Value = 2.4
$ws.Range("L271").Value = 3
$ws.Range("M271").Value = 2.9
$ws.Range("N271").Value = 2.3
$ws.Range("O271").Value = 3
$ws.Range("P271").Value = 3.1
$ws.Range("Q271").Value = -0.25
$ws.Range("R271").Value = 1.975
$ws.Range("S271").Value = 1.875
$ws.Range("T271").Value = 2.25
$ws.Range("U271").Value = 1.875
$ws.Range("V271").Value = 1.975
$ws.Range("W271").Value = 1.3
$ws.Range("X271").Value = -1
$ws.Range("Y271").Value = -1
$ws.Range("Z271").Value = 0.9750000000000001
$ws.Range("AA271").Value = -1
$ws.Range("AB271").Value = 0.875
$ws.Range("AC271").Value = -1
$ws.Range("A235").Copy()
$ws.Range("A271").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E271").PasteSpecial(-4122)

# Row 272
$ws.Range("A272").Value = 270
$ws.Range("B272").Value = 7951789
$ws.Range("C272").Value = "Romania Liga I"
$ws.Range("D272").Value = "Romania Liga I"
$ws.Range("E272").Value = 45396.51041666666
$ws.Range("F272").Value = "Dinamo Bucharest"
$ws.Range("G272").Value = "CSM Politehnica Iasi"
$ws.Range("H272").Value = 1
$ws.Range("I272").Value = 0
$ws.Range("J272").Value = "H"
$ws.Range("K272").Value = 2.25
$ws.Range("L272").Value = 3.1
$ws.Range("M272").Value = 3.1
$ws.Range("N272").Value = 2.1
$ws.Range("O272").Value = 3.1
$ws.Range("P272").Value = 3.3
$ws.Range("Q272").Value = -0.25
$ws.Range("R272").Value = 1.875
$ws.Range("S272").Value = 1.975
$ws.Range("T272").Value = 2.25
$ws.Range("U272").Value = 2
$ws.Range("V272").Value = 1.85
$ws.Range("W272").Value = 1.1
$ws.Range("X272").Value = -1
$ws.Range("Y272").Value = -1
$ws.Range("Z272").Value = 0.875
$ws.Range("AA272").Value = -1
$ws.Range("AB272").Value = -1
$ws.Range("AC272").Value = 0.8500000000000001
$ws.Range("A235").Copy()
$ws.Range("A272").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E272").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Edit complete. UsedRange:" $($ws.UsedRange.Address())